# Converts an EMU (English Metric Unit) value to points for the COM
# Left/Top/Width/Height properties. A tiny epsilon is added to counteract
# truncation in the host's point->EMU re-serialization so the round trip
# lands on the exact EMU value.
function EMU($v) {
    return ($v / 12700.0) + 0.00002
}

$p = $ppt.ActivePresentation

# The new slide 6 is a duplicate of slide 5 ("BaB" board-game slide) with
# its "Bad move!" feedback textbox + its backing rectangle swapped out for
# a second "BaB" label rectangle (used when there is no feedback message).
$s5 = $p.Slides.Item(5)
$dupSlides = $s5.Duplicate()
$s6 = $dupSlides.Item(1)

# Remove "Rectangle 18" (blank bg1 rectangle behind the feedback message)
# and "TextBox 5" (the "Bad move! ..." feedback text).
$s6.Shapes.Item(18).Delete()
$s6.Shapes.Item(17).Delete()

# Duplicate the existing "Rectangle 17" ("BaB" label) shape to create the
# new "Rectangle 20" shape, then move/resize it into place.
$srcRect = $s6.Shapes.Item(16)
$newRectColl = $srcRect.Duplicate()
$newRect = $newRectColl.Item(1)
$newRect.Name = "Rectangle 20"
$newRect.Left = EMU(4037372)
$newRect.Top = EMU(3672945)
$newRect.Width = EMU(4768990)
$newRect.Height = EMU(746655)
